$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.427.81'
$ws.Range('E2').Value = '  +2.41%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.400.61'
$ws.Range('E3').Value = '  +1.42%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '561.09'
$ws.Range('E5').Value = '  +2.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '176.05'
$ws.Range('E6').Value = '  +1.97%  '
$ws.Range('E7').Value = '  +2.68%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.393.68'
$ws.Range('E8').Value = '  +1.54%  '
$ws.Range('E9').Value = '  -0.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.170'
$ws.Range('E10').Value = '  +12.05%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.633'
$ws.Range('E11').Value = '  +3.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.93'
$ws.Range('E12').Value = '  +1.93%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000279'
$ws.Range('E13').Value = '  +5.17%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.16'
$ws.Range('E14').Value = '  +2.63%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.947.04'
$ws.Range('E15').Value = '  +2.08%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.35'
$ws.Range('E16').Value = '  +2.18%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.404.07'
$ws.Range('E17').Value = '  +1.96%  '
$ws.Range('E18').Value = '  +1.60%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '65.475.54'
$ws.Range('E19').Value = '  +2.68%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.89'
$ws.Range('E20').Value = '  +1.32%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.995'
$ws.Range('E21').Value = '  +1.77%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '473.50'
$ws.Range('E22').Value = '  +14.67%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.07'
$ws.Range('E23').Value = '  +17.31%  '
$ws.Range('E24').Value = '  +2.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '86.82'
$ws.Range('E25').Value = '  +4.46%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '13.45'
$ws.Range('E26').Value = '  -2.64%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.92'
$ws.Range('E27').Value = '  +3.24%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.89'
$ws.Range('E28').Value = '  +5.92%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.87'
$ws.Range('E29').Value = '  +2.80%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.17'
$ws.Range('E30').Value = '  +7.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.73'
$ws.Range('E31').Value = '  +4.79%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.55'
$ws.Range('E32').Value = '  +1.77%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '62.52'
$ws.Range('E33').Value = '  +7.77%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '575.63'
$ws.Range('E34').Value = '  +0.70%  '
$ws.Range('E35').Value = '  +1.67%  '
$ws.Range('E36').Value = '  -0.21%  '
$ws.Range('B37').Value = 'Stacks'
$ws.Range('C37').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.54'
$ws.Range('E37').Value = '  +4.12%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.140'
$ws.Range('E38').Value = '  -5.01%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.85'
$ws.Range('E39').Value = '  +1.71%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0₃0757'
$ws.Range('E40').Value = '  +2.20%  '
$ws.Range('E41').Value = '  +1.64%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.091.90'
$ws.Range('E42').Value = '  -1.88%  '
$ws.Range('E43').Value = '  +0.37%  '
$ws.Range('E44').Value = '  +1.40%  '
$ws.Range('E45').Value = '  +4.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.49'
$ws.Range('E46').Value = '  +2.98%  '
$ws.Range('E47').Value = '  +5.58%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.18'
$ws.Range('E48').Value = '  -2.83%  '
$ws.Range('E49').Value = '  +0.49%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '137.15'
$ws.Range('E50').Value = '  +3.77%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.33'
$ws.Range('E51').Value = '  +2.89%  '
